$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price/Volume columns remain plain text (they store strings like
# "30.265.69" or "  -0.22%  " which Excel would otherwise auto-convert to
# numbers/percentages).
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("E2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "30.265.69"
$ws.Range("E2").Value = "  -0.22%  "
$ws.Range("D3").Value = "1.928.84"
$ws.Range("E3").Value = "  -0.13%  "
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  +0.26%  "
$ws.Range("D5").Value = "249.10"
$ws.Range("E5").Value = "  +0.08%  "
$ws.Range("D6").Value = "0.7129"
$ws.Range("E6").Value = "  -1.05%  "
$ws.Range("E7").Value = "  +0.18%  "
$ws.Range("D8").Value = "0.3211"
$ws.Range("E8").Value = "  -1.81%  "
$ws.Range("D9").Value = "27.08"
$ws.Range("E9").Value = "  -0.45%  "
$ws.Range("D10").Value = "0.07089"
$ws.Range("E10").Value = "  +4.25%  "
$ws.Range("D11").Value = "0.7902"
$ws.Range("E11").Value = "  -1.64%  "
$ws.Range("D12").Value = "0.08028"
$ws.Range("E12").Value = "  -0.30%  "
$ws.Range("D13").Value = "1.927.78"
$ws.Range("E13").Value = "  -0.23%  "
$ws.Range("D14").Value = "5.365"
$ws.Range("E14").Value = "  -0.75%  "
$ws.Range("D15").Value = "94.52"
$ws.Range("E15").Value = "  -0.22%  "
$ws.Range("D16").Value = "14.58"
$ws.Range("E16").Value = "  +0.88%  "
$ws.Range("D17").Value = "30.279.28"
$ws.Range("E17").Value = "  -0.15%  "
$ws.Range("D18").Value = "256.25"
$ws.Range("E18").Value = "  +1.62%  "
$ws.Range("D19").Value = "0.000008047"
$ws.Range("E19").Value = "  +0.85%  "
$ws.Range("D20").Value = "5.749"
$ws.Range("E20").Value = "  -1.16%  "
$ws.Range("D21").Value = "2.181.88"
$ws.Range("E21").Value = "  -0.21%  "
$ws.Range("E22").Value = "  +0.18%  "
$ws.Range("E23").Value = "  +0.08%  "
$ws.Range("D24").Value = "6.796"
$ws.Range("E24").Value = "  -0.96%  "
$ws.Range("D25").Value = "9.530"
$ws.Range("E25").Value = "  -1.33%  "
$ws.Range("D26").Value = "165.67"
$ws.Range("E26").Value = "  +3.91%  "
$ws.Range("D27").Value = "19.15"
$ws.Range("E27").Value = "  +0.54%  "
$ws.Range("D28").Value = "2.273"
$ws.Range("E28").Value = "  -4.66%  "
$ws.Range("D29").Value = "0.1271"
$ws.Range("E29").Value = "  -4.40%  "
$ws.Range("E30").Value = "  +1.27%  "
$ws.Range("D31").Value = "1.529"
$ws.Range("E31").Value = "  -1.84%  "
$ws.Range("D32").Value = "4.389"
$ws.Range("E32").Value = "  -0.04%  "
$ws.Range("D33").Value = "4.117"
$ws.Range("E33").Value = "  -1.56%  "
$ws.Range("D34").Value = "0.05152"
$ws.Range("E34").Value = "  +1.87%  "
$ws.Range("E35").Value = "  +3.16%  "
$ws.Range("E36").Value = "  +0.76%  "
$ws.Range("D37").Value = "2.768"
$ws.Range("D38").Value = "0.01952"
$ws.Range("E38").Value = "  -0.76%  "
$ws.Range("E39").Value = "  -0.66%  "
$ws.Range("D40").Value = "77.52"
$ws.Range("E40").Value = "  -2.15%  "
$ws.Range("D41").Value = "6.343"
$ws.Range("E41").Value = "  -3.82%  "
$ws.Range("D42").Value = "0.4471"
$ws.Range("E42").Value = "  +0.44%  "
$ws.Range("D43").Value = "1.980"
$ws.Range("E43").Value = "  -0.48%  "
$ws.Range("D44").Value = "0.8426"
$ws.Range("E44").Value = "  +1.06%  "
$ws.Range("D45").Value = "1.002"
$ws.Range("E45").Value = "  +0.10%  "
$ws.Range("D46").Value = "100.93"
$ws.Range("E46").Value = "  -1.09%  "
$ws.Range("D47").Value = "9.705"
$ws.Range("E47").Value = "  -0.54%  "
$ws.Range("D48").Value = "7.412"
$ws.Range("E48").Value = "  +2.00%  "
$ws.Range("D49").Value = "36.39"
$ws.Range("E49").Value = "  +0.12%  "
$ws.Range("D50").Value = "0.06110"
$ws.Range("E50").Value = "  +2.92%  "
$ws.Range("D51").Value = "0.4156"
$ws.Range("E51").Value = "  +2.39%  "
